# Generate Report for Archive
#
# Two logical changes are applied:
#  1. The localization status text "Ready for handoff" -> "In Translation"
#     everywhere it appears (Overview!E2:F2, zh-cn!C2, de-de!C2). All cells
#     sharing that string are updated together so the old string is fully
#     retired from the shared-string table rather than just duplicated.
#  2. The two "handoff/handback datetime" columns on each sheet are made
#     narrower (they were sized the same as the long filename columns;
#     now they match the narrower date columns instead).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# --- 2. Narrow the datetime columns on each sheet ---
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
